# Updates cryptos list: prices (col D) and 1h volume/change % (col E)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "26.361.25"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.592.86"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.43%  "
$c = $ws.Range("D5")
$c.Value = "'210.23"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -0.41%  "
$c = $ws.Range("D8")
$c.Value = "'0.0612"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.816.22"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "1.532.67"
$ws.Range("E14").Value = "  -4.57%  "
$ws.Range("E15").Value = "  -1.61%  "
$c = $ws.Range("D16")
$c.Value = "'64.68"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "26.349.85"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.94%  "
$c = $ws.Range("D19")
$c.Value = "'7.49"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.08%  "
$c = $ws.Range("D20")
$c.Value = "'211.76"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("E24").Value = "  -1.42%  "
$c = $ws.Range("D25")
$c.Value = "'145.39"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -0.39%  "
$c = $ws.Range("D27")
$c.Value = "'7.06"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  -1.08%  "
$c = $ws.Range("D29")
$c.Value = "'15.29"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.63%  "
$c = $ws.Range("D30")
$c.Value = "'0.0504"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "1.301.00"
$ws.Range("E34").Value = "  +1.38%  "
$c = $ws.Range("D35")
$c.Value = "'0.613"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  -0.40%  "
$c = $ws.Range("D39")
$c.Value = "'1.11"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -13.60%  "
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  +2.93%  "
$c = $ws.Range("D43")
$c.Value = "'62.73"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("E44").Value = "  -5.11%  "
$c = $ws.Range("D45")
$c.Value = "'0.762"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").Value = "1.728.65"
$ws.Range("E46").Value = "  -0.61%  "
$c = $ws.Range("D47")
$c.Value = "'88.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("E49").Value = "  +5.84%  "
$c = $ws.Range("D50")
$c.Value = "'0.0986"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.20%  "
